# Update "想去人数" (attendance/interest count) figures for three events
# that appear on both the "展览" (Exhibitions) sheet and the "全部类型"
# (All types) sheet.
#   熊喵M动漫嘉年华【免费】      1243 -> 1247
#   第二届北极光动漫展            2734 -> 2736
#   万圣漫控嘉年华10              244  -> 245

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1247
$ws1.Range("F4").Value = 2736
$ws1.Range("F5").Value = 245

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1247
$ws4.Range("F6").Value = 2736
$ws4.Range("F8").Value = 245
